$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76
$ws.Cells.Item(76, 1).Value = 'q'
$ws.Cells.Item(76, 2).Value = 'q'
$ws.Cells.Item(76, 3).Value = 'Número inválido'
$ws.Cells.Item(76, 4).Value = 'qq'
$ws.Cells.Item(76, 5).Value = 'q'
$ws.Cells.Item(76, 6).Value = 'q'
$ws.Cells.Item(76, 8).Value = 'q'

# Row 77
$ws.Cells.Item(77, 1).Value = 'q'
$ws.Cells.Item(77, 2).Value = 'q'
$ws.Cells.Item(77, 3).Value = 'Número inválido'
$ws.Cells.Item(77, 4).Value = 'q'
$ws.Cells.Item(77, 5).Value = 'q'
$ws.Cells.Item(77, 6).Value = 'qq'
$ws.Cells.Item(77, 7).Value = 'TAX ID'
$ws.Cells.Item(77, 8).Value = 'q'

# Row 78
$ws.Cells.Item(78, 1).Value = 'q'
$ws.Cells.Item(78, 2).Value = 'q'
$ws.Cells.Item(78, 3).Value = 'Número inválido'
$ws.Cells.Item(78, 4).Value = 'q'
$ws.Cells.Item(78, 5).Value = 'q'
$ws.Cells.Item(78, 6).Value = 'q'
$ws.Cells.Item(78, 8).Value = 'q'

# Row 79
$ws.Cells.Item(79, 1).Value = 'q'
$ws.Cells.Item(79, 2).Value = 'q'
$ws.Cells.Item(79, 3).Value = 'Número inválido'
$ws.Cells.Item(79, 4).Value = 'q'
$ws.Cells.Item(79, 5).Value = 'q'
$ws.Cells.Item(79, 6).Value = 'q'
$ws.Cells.Item(79, 7).Value = 'TAX ID'
$ws.Cells.Item(79, 8).Value = 'q'

# Row 80
$ws.Cells.Item(80, 1).Value = 'q'
$ws.Cells.Item(80, 2).Value = 'q'
$ws.Cells.Item(80, 3).Value = 'Número inválido'
$ws.Cells.Item(80, 4).Value = 'q'
$ws.Cells.Item(80, 5).Value = 'q'
$ws.Cells.Item(80, 6).Value = 'q'
$ws.Cells.Item(80, 8).Value = 'q'

# Row 81
$ws.Cells.Item(81, 1).Value = 'q'
$ws.Cells.Item(81, 2).Value = 'q'
$ws.Cells.Item(81, 3).Value = 'Número inválido'
$ws.Cells.Item(81, 4).Value = 'q'
$ws.Cells.Item(81, 5).Value = 'q'
$ws.Cells.Item(81, 6).Value = 'q'
$ws.Cells.Item(81, 7).Value = 'TAX ID'
$ws.Cells.Item(81, 8).Value = 'q'

# Row 82
$ws.Cells.Item(82, 1).Value = 'q'
$ws.Cells.Item(82, 2).Value = 'q'
$ws.Cells.Item(82, 3).Value = 'Número inválido'
$ws.Cells.Item(82, 4).Value = 'q'
$ws.Cells.Item(82, 5).Value = 'q'
$ws.Cells.Item(82, 6).Value = 'q'
$ws.Cells.Item(82, 8).Value = 'q'

# Row 83
$ws.Cells.Item(83, 1).Value = 'q'
$ws.Cells.Item(83, 2).Value = 'q'
$ws.Cells.Item(83, 3).Value = 'Número inválido'
$ws.Cells.Item(83, 4).Value = 'q'
$ws.Cells.Item(83, 5).Value = 'q'
$ws.Cells.Item(83, 6).Value = 'q'
$ws.Cells.Item(83, 7).Value = 'q'
$ws.Cells.Item(83, 8).Value = 'q'

# Row 84
$ws.Cells.Item(84, 1).Value = 'q'
$ws.Cells.Item(84, 2).Value = 'q'
$ws.Cells.Item(84, 3).Value = 'Número inválido'
$ws.Cells.Item(84, 4).Value = 'q'
$ws.Cells.Item(84, 5).Value = 'q'
$ws.Cells.Item(84, 6).Value = 'q'
$ws.Cells.Item(84, 8).Value = 'q'

# Row 85
$ws.Cells.Item(85, 1).Value = 'q'
$ws.Cells.Item(85, 2).Value = 'q'
$ws.Cells.Item(85, 3).Value = 'Número inválido'
$ws.Cells.Item(85, 4).Value = 'q'
$ws.Cells.Item(85, 5).Value = 'q'
$ws.Cells.Item(85, 6).Value = 'q'
$ws.Cells.Item(85, 8).Value = 'q'

# Row 86
$ws.Cells.Item(86, 1).Value = 'q'
$ws.Cells.Item(86, 2).Value = 'q'
$ws.Cells.Item(86, 3).Value = 'Número inválido'
$ws.Cells.Item(86, 4).Value = 'q'
$ws.Cells.Item(86, 5).Value = 'q'
$ws.Cells.Item(86, 6).Value = 'q'
$ws.Cells.Item(86, 7).Value = 'qqq'
$ws.Cells.Item(86, 8).Value = 'q'

# Row 87
$ws.Cells.Item(87, 1).Value = 'q'
$ws.Cells.Item(87, 2).Value = 'q'
$ws.Cells.Item(87, 3).Value = 'Número inválido'
$ws.Cells.Item(87, 4).Value = 'q'
$ws.Cells.Item(87, 5).Value = 'q'
$ws.Cells.Item(87, 6).Value = 'q'
$ws.Cells.Item(87, 8).Value = 'q'

# Row 88
$ws.Cells.Item(88, 1).Value = 'h'
$ws.Cells.Item(88, 2).Value = 'h'
$ws.Cells.Item(88, 3).Value = 'Número inválido'
$ws.Cells.Item(88, 4).Value = 'h'
$ws.Cells.Item(88, 5).Value = 'h'
$ws.Cells.Item(88, 6).Value = 'h'
$ws.Cells.Item(88, 7).Value = 'h'
$ws.Cells.Item(88, 8).Value = 'h'

# Row 89
$ws.Cells.Item(89, 1).Value = 'w'
$ws.Cells.Item(89, 2).Value = 'w'
$ws.Cells.Item(89, 3).Value = 'Número inválido'
$ws.Cells.Item(89, 4).Value = 'w'
$ws.Cells.Item(89, 5).Value = 'w'
$ws.Cells.Item(89, 6).Value = 'w'
$ws.Cells.Item(89, 8).Value = 'w'

# Row 90
$ws.Cells.Item(90, 1).Value = 'q'
$ws.Cells.Item(90, 2).Value = 'q'
$ws.Cells.Item(90, 3).Value = 'Número inválido'
$ws.Cells.Item(90, 4).Value = 'q'
$ws.Cells.Item(90, 5).Value = 'q'
$ws.Cells.Item(90, 6).Value = 'q'
$ws.Cells.Item(90, 8).Value = 'q'

# Row 91
$ws.Cells.Item(91, 1).Value = 's'
$ws.Cells.Item(91, 2).Value = 's'
$ws.Cells.Item(91, 3).Value = 'Número inválido'
$ws.Cells.Item(91, 4).Value = 's'
$ws.Cells.Item(91, 5).Value = 's'
$ws.Cells.Item(91, 6).Value = 's'
$ws.Cells.Item(91, 7).Value = 's'
$ws.Cells.Item(91, 8).Value = 's'

# Row 92
$ws.Cells.Item(92, 1).Value = 'd'
$ws.Cells.Item(92, 2).Value = 'd'
$ws.Cells.Item(92, 3).Value = 'Número inválido'
$ws.Cells.Item(92, 4).Value = 'd'
$ws.Cells.Item(92, 5).Value = 'd'
$ws.Cells.Item(92, 6).Value = 'd'
$ws.Cells.Item(92, 8).Value = 'd'

# Row 93
$ws.Cells.Item(93, 1).Value = 'ss'
$ws.Cells.Item(93, 2).Value = 's'
$ws.Cells.Item(93, 3).Value = 'Número inválido'
$ws.Cells.Item(93, 4).Value = 's'
$ws.Cells.Item(93, 5).Value = 's'
$ws.Cells.Item(93, 6).Value = 's'
$ws.Cells.Item(93, 7).Value = 's'
$ws.Cells.Item(93, 8).Value = 's'

# Row 94
$ws.Cells.Item(94, 1).Value = 's'
$ws.Cells.Item(94, 2).Value = 's'
$ws.Cells.Item(94, 3).Value = 'Número inválido'
$ws.Cells.Item(94, 4).Value = 's'
$ws.Cells.Item(94, 5).Value = 's'
$ws.Cells.Item(94, 6).Value = 's'
$ws.Cells.Item(94, 8).Value = 's'

# Row 95
$ws.Cells.Item(95, 1).Value = 'w'
$ws.Cells.Item(95, 2).Value = 'w'
$ws.Cells.Item(95, 3).Value = 'Número inválido'
$ws.Cells.Item(95, 4).Value = 'w'
$ws.Cells.Item(95, 5).Value = 'w'
$ws.Cells.Item(95, 6).Value = 'w'
$ws.Cells.Item(95, 7).Value = 'w'
$ws.Cells.Item(95, 8).Value = 'w'

# Row 96
$ws.Cells.Item(96, 1).Value = 'd'
$ws.Cells.Item(96, 2).Value = 'd'
$ws.Cells.Item(96, 3).Value = 'Número inválido'
$ws.Cells.Item(96, 4).Value = 'd'
$ws.Cells.Item(96, 5).Value = 'd'
$ws.Cells.Item(96, 6).Value = 'd'
$ws.Cells.Item(96, 8).Value = 'd'

# Row 97
$ws.Cells.Item(97, 1).Value = 'e'
$ws.Cells.Item(97, 2).Value = 'e'
$ws.Cells.Item(97, 3).Value = 'Número inválido'
$ws.Cells.Item(97, 4).Value = 'e'
$ws.Cells.Item(97, 5).Value = 'e'
$ws.Cells.Item(97, 6).Value = 'e'
$ws.Cells.Item(97, 7).Value = 'e'
$ws.Cells.Item(97, 8).Value = 'e'
